# Update cryptos list: refresh Price / Volume(1h) figures, and swap a handful
# of adjacent coin rows (Stacks<->RenderToken, USDe<->WhiteBITCoin,
# Mantle<->InjectiveProtocol). Price-looking strings are written with a
# leading "'" so Excel stores them as text (matching the original
# inlineStr cells) instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.959.75'
$ws.Range("E2").Value = '''  +1.31%  '

$ws.Range("D3").Value = '''2.621.76'
$ws.Range("E3").Value = '''  +1.20%  '

$ws.Range("E4").Value = '''  +0.39%  '

$ws.Range("D5").Value = '''596.13'
$ws.Range("E5").Value = '''  +0.59%  '

$ws.Range("D6").Value = '''154.88'
$ws.Range("E6").Value = '''  -0.28%  '

$ws.Range("E7").Value = '''  +0.09%  '

$ws.Range("E8").Value = '''  +0.72%  '

$ws.Range("D9").Value = '''2.618.67'
$ws.Range("E9").Value = '''  +1.07%  '

$ws.Range("D10").Value = '''0.127'
$ws.Range("E10").Value = '''  +9.30%  '

$ws.Range("E11").Value = '''  +1.00%  '

$ws.Range("E12").Value = '''  +0.29%  '

$ws.Range("E13").Value = '''  -1.52%  '

$ws.Range("D14").Value = '''27.76'
$ws.Range("E14").Value = '''  -2.30%  '

$ws.Range("E15").Value = '''  +3.20%  '

$ws.Range("D16").Value = '''3.102.03'
$ws.Range("E16").Value = '''  +1.54%  '

$ws.Range("D17").Value = '''67.769.60'
$ws.Range("E17").Value = '''  +1.53%  '

$ws.Range("D18").Value = '''2.625.16'

$ws.Range("D19").Value = '''367.33'
$ws.Range("E19").Value = '''  +3.14%  '

$ws.Range("D20").Value = '''11.18'
$ws.Range("E20").Value = '''  -1.75%  '

$ws.Range("D21").Value = '''7.64'
$ws.Range("E21").Value = '''  -2.51%  '

$ws.Range("E22").Value = '''  -0.83%  '

$ws.Range("D23").Value = '''2.05'
$ws.Range("E23").Value = '''  -1.14%  '

$ws.Range("E24").Value = '''  +0.03%  '

$ws.Range("D25").Value = '''9.90'
$ws.Range("E25").Value = '''  -6.51%  '

$ws.Range("D26").Value = '''67.58'
$ws.Range("E26").Value = '''  +0.77%  '

$ws.Range("E27").Value = '''  +0.78%  '

$ws.Range("D28").Value = '''2.730.84'
$ws.Range("E28").Value = '''  +0.67%  '

$ws.Range("D29").Value = '''579.93'
$ws.Range("E29").Value = '''  -4.78%  '

$ws.Range("D30").Value = '''1.04'
$ws.Range("E30").Value = '''  +4.40%  '

$ws.Range("D31").Value = '''1.43'
$ws.Range("E31").Value = '''  -2.44%  '

$ws.Range("E32").Value = '''  -1.35%  '

$ws.Range("E33").Value = '''  +0.11%  '

$ws.Range("E34").Value = '''  -1.94%  '

$ws.Range("D35").Value = '''1.00'

$ws.Range("D36").Value = '''1.52'
$ws.Range("E36").Value = '''  -4.03%  '

$ws.Range("D37").Value = '''4.93'
$ws.Range("E37").Value = '''  -2.77%  '

$ws.Range("D38").Value = '''158.53'
$ws.Range("E38").Value = '''  +2.79%  '

$ws.Range("D39").Value = '''19.36'
$ws.Range("E39").Value = '''  +0.13%  '

$ws.Range("D40").Value = '''0.370'
$ws.Range("E40").Value = '''  -0.41%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").Value = '''5.33'
$ws.Range("E41").Value = '''  -3.30%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''1.85'
$ws.Range("E42").Value = '''  +1.62%  '

$ws.Range("E43").Value = '''  -4.10%  '

$ws.Range("D44").Value = '''41.19'
$ws.Range("E44").Value = '''  -0.87%  '

$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").Value = '''16.43'
$ws.Range("E45").Value = '''  -0.12%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '''1.00'
$ws.Range("E46").Value = '''  +0.03%  '

$ws.Range("D47").Value = '''156.54'
$ws.Range("E47").Value = '''  -0.26%  '

$ws.Range("D48").Value = '''0.0₆0286'
$ws.Range("E48").Value = '''  -7.77%  '

$ws.Range("E49").Value = '''  -0.64%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''20.99'
$ws.Range("E50").Value = '''  -2.64%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.625'
$ws.Range("E51").Value = '''  +1.81%  '
